# Timesheet updated by Kirubaharan
#
# Applies the edits described by the commit diff:
#  - Sheet "6-4-22 " (Worksheets(2)): view scrolled back to top, selection -> C15
#  - Sheet "7-4-22"  (Worksheets(3)): several cell edits (rows 5, 7, 10, 12),
#    new dimension/selection, selection -> F12

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "6-4-22 " - only the view/selection changed (topLeftCell cleared,
# active cell moves from C13 to C15).
# ---------------------------------------------------------------------------
$wsB = $wb.Worksheets.Item(2)
$wsB.Activate()
$wsB.Range("C15").Select()

# ---------------------------------------------------------------------------
# Sheet "7-4-22" - cell content updates.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item(3)
$ws.Activate()

# Row 5 (Iswarya)
$ws.Cells.Item(5, 2).ClearContents()
$ws.Cells.Item(5, 3).Value = "UI designing on trainer - landing page,feedback page"
$ws.Cells.Item(5, 5).Clear()
$ws.Cells.Item(5, 8).Value = "Trainer prototype completed,session about abstract pattern,Prototype disscussion with rafi,team discussion about modification on prototype"

# Row 7 (Kirubaharan)
$ws.Cells.Item(7, 2).Value = "Prototype discussion with team "
$ws.Cells.Item(7, 3).Value = "Assumption, Acceptance Criteria and constraints for Trainee stories"
$ws.Cells.Item(7, 4).Value = "Prototype design"
$ws.Cells.Item(7, 5).Value = " Discussion with team members about UI prototype - 40 mins, Worked on Acceptance criteria ,assumption and constraints for Trainee - 2 hours,Meeting with Rafi - 1 hour,  "
$ws.Cells.Item(7, 6).Value = "Design pattern session - 40 mins ,softskill session- 50 mins,Others - 1 hour 20 mins"

# Row 10 (Sruthi)
$ws.Cells.Item(10, 3).Value = "Trainer - constraints, acceptance criteria, assumptions"
$ws.Cells.Item(10, 4).Clear()
$ws.Cells.Item(10, 5).Value = "Wor,session about abstract pattern,Prototype disscussion with rafi,team discussion about modification on prototype"

# Row 12 (Arul)
$ws.Cells.Item(12, 1).Value = "Arul "
$ws.Cells.Item(12, 2).Value = "Refining trainee user stories"
$ws.Cells.Item(12, 3).Value = "Assumption, Acceptance Criteria and constraints for Trainee stories"
$ws.Cells.Item(12, 4).Value = "Prototype design"
$ws.Cells.Item(12, 5).Value = "Meeting with client - 45mins,I've discussed about UI protyping with team members - 45mins, Done acceptance criteria for Trainee user stories - 120 mins"
$ws.Cells.Item(12, 6).Value = "Softskill Session-45 mins, Design pattern session ,others 90mins"

# Final view state for sheet "7-4-22": it stays the active/visible tab with
# the selection parked on F12.
$ws.Activate()
$ws.Range("F12").Select()
